# Append run: refresh the "取得日時" (fetched-at) timestamp on every data row,
# and pick up the re-ordered (swapped) title/URL pair for the two SRE listings
# (rows 7 & 8) that the scraper returned this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-16 12:42:00"

# 1) Update the "取得日時" column (A) for every data row (2-10) to the new run time.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 2) Rows 7 and 8 swapped their title (B) and URL (F) values between runs.
$b7 = $ws.Range("B7").Value2
$f7 = $ws.Range("F7").Value2
$b8 = $ws.Range("B8").Value2
$f8 = $ws.Range("F8").Value2

$ws.Range("B7").Value = $b8
$ws.Range("F7").Value = $f8
$ws.Range("B8").Value = $b7
$ws.Range("F8").Value = $f7
